$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; existing rows 17..127 shift down to 18..128
$ws.Rows(17).Insert()

# Populate the newly inserted row 17 with the new observation
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = "2021-12-07"
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112009
$ws.Range("G17").Value = "Acelga"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 3500
$ws.Range("L17").Value = 3500
$ws.Range("M17").Value = 3500
$ws.Range("N17").Value = '$/docena de atados (4 kilos)'
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 875
$ws.Range("Q17").Value = 4
$ws.Range("R17").Value = "Hortaliza"
